$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LFC")

$ws.Range("D8").Value = 59803000
$ws.Range("E8").Value = 36844300
$ws.Range("F8").Value = 59045200
$ws.Range("G8").Value = 30384100
$ws.Range("H8").Value = 50267400
$ws.Range("I8").Value = 26181800
$ws.Range("J8").Value = 49176400

$ws.Range("D14").Value = 188200
$ws.Range("E14").Value = 158800
$ws.Range("F14").Value = 250400
$ws.Range("G14").Value = 250400
$ws.Range("H14").Value = 143800
$ws.Range("I14").Value = 42000
$ws.Range("J14").Value = 5600

$ws.Range("D17").Value = 56918100
$ws.Range("E17").Value = 33229300
$ws.Range("F17").Value = 56853100
$ws.Range("G17").Value = 29074100
$ws.Range("H17").Value = 48200500
$ws.Range("I17").Value = 25162000
$ws.Range("J17").Value = 43031300

$ws.Range("D18").Value = 2884900
$ws.Range("E18").Value = 3615000
$ws.Range("F18").Value = 2192200
$ws.Range("G18").Value = 1310000
$ws.Range("H18").Value = 2066900
$ws.Range("I18").Value = 1019700
$ws.Range("J18").Value = 6145100

$ws.Range("D20").Value = 298000
$ws.Range("E20").Value = 205400
$ws.Range("F20").Value = 171900
$ws.Range("G20").Value = 265200
$ws.Range("H20").Value = -103700
$ws.Range("I20").Value = -323200
$ws.Range("J20").Value = -24900

$ws.Range("D21").Value = 3370200
$ws.Range("E21").Value = 3828700
$ws.Range("F21").Value = 2526100
$ws.Range("G21").Value = 1576900
$ws.Range("H21").Value = 2116900
$ws.Range("I21").Value = 694400
$ws.Range("J21").Value = 6272300

$ws.Range("D23").Value = 3182900
$ws.Range("E23").Value = 3820400
$ws.Range("F23").Value = 2364000
$ws.Range("G23").Value = 1575200
$ws.Range("H23").Value = 1963200
$ws.Range("I23").Value = 696500
$ws.Range("J23").Value = 6120100

$ws.Range("D24").Value = 704100
$ws.Range("E24").Value = 813600
$ws.Range("F24").Value = 510100
$ws.Range("G24").Value = 248700
$ws.Range("H24").Value = 383000
$ws.Range("I24").Value = 184000
$ws.Range("J24").Value = 1410500

$ws.Range("D26").Value = 2478900
$ws.Range("E26").Value = 3006800
$ws.Range("F26").Value = 1853900
$ws.Range("G26").Value = 1326500
$ws.Range("H26").Value = 1580100
$ws.Range("I26").Value = 512500
$ws.Range("J26").Value = 4709600

$ws.Range("D27").Value = 2437300
$ws.Range("E27").Value = 2969800
$ws.Range("F27").Value = 1816800
$ws.Range("G27").Value = 1295900
$ws.Range("H27").Value = 1542700
$ws.Range("I27").Value = 476400
$ws.Range("J27").Value = 4673300

$ws.Range("D32").Value = -298000
$ws.Range("E32").Value = -205400
$ws.Range("F32").Value = -171900
$ws.Range("G32").Value = -265200
$ws.Range("H32").Value = 103700
$ws.Range("I32").Value = 323200
$ws.Range("J32").Value = 24900

$ws.Range("D33").Value = 2437300
$ws.Range("E33").Value = 2969800
$ws.Range("F33").Value = 1816800
$ws.Range("G33").Value = 1295900
$ws.Range("H33").Value = 1542700
$ws.Range("I33").Value = 476400
$ws.Range("J33").Value = 4673300

$ws.Range("D35").Value = 2437300
$ws.Range("E35").Value = 2969800
$ws.Range("F35").Value = 1816800
$ws.Range("G35").Value = 1295900
$ws.Range("H35").Value = 1542700
$ws.Range("I35").Value = 476400
$ws.Range("J35").Value = 4673300

$ws.Range("D41").Value = 8603800
$ws.Range("E41").Value = 7210600
$ws.Range("F41").Value = 6907400
$ws.Range("G41").Value = 9552300
$ws.Range("H41").Value = 10658700
$ws.Range("I41").Value = 11002400
$ws.Range("J41").Value = 10970200

$ws.Range("D43").Value = 5605300
$ws.Range("E43").Value = 2095700
$ws.Range("F43").Value = 4566100
$ws.Range("G43").Value = 1991800
$ws.Range("H43").Value = 4094900
$ws.Range("I43").Value = 1768000
$ws.Range("J43").Value = 3400100

$ws.Range("D47").Value = 417285000
$ws.Range("E47").Value = 399985900
$ws.Range("F47").Value = 395506600
$ws.Range("G47").Value = 371197300
$ws.Range("H47").Value = 353073200
$ws.Range("I47").Value = 334383900
$ws.Range("J47").Value = 323344700

$ws.Range("D48").Value = 7554500
$ws.Range("E48").Value = 6792900
$ws.Range("F48").Value = 5883100
$ws.Range("G48").Value = 4686800
$ws.Range("H48").Value = 4228900
$ws.Range("I48").Value = 4186800
$ws.Range("J48").Value = 3925400

$ws.Range("D54").Value = 451638000
$ws.Range("E54").Value = 430031500
$ws.Range("F54").Value = 426692400
$ws.Range("G54").Value = 400254500
$ws.Range("H54").Value = 384240300
$ws.Range("I54").Value = 363354400
$ws.Range("J54").Value = 353149800

$ws.Range("D59").Value = 381997800
$ws.Range("E59").Value = 362087300
$ws.Range("F59").Value = 360110900
$ws.Range("G59").Value = 332057400
$ws.Range("H59").Value = 314635300
$ws.Range("I59").Value = 293381600
$ws.Range("J59").Value = 285046800

$ws.Range("D61").Value = 2902500
$ws.Range("E61").Value = 2789200
$ws.Range("F61").Value = 3879000
$ws.Range("G61").Value = 8039100
$ws.Range("H61").Value = 10455200
$ws.Range("I61").Value = 10483200
$ws.Range("J61").Value = 10484000

$ws.Range("D62").Value = 463200
$ws.Range("E62").Value = 722900
$ws.Range("F62").Value = 1398900
$ws.Range("G62").Value = 1152800
$ws.Range("H62").Value = 1706000
$ws.Range("I62").Value = 2516000
$ws.Range("J62").Value = 2839200

$ws.Range("D66").Value = 403700400
$ws.Range("E66").Value = 382401800
$ws.Range("F66").Value = 380893800
$ws.Range("G66").Value = 355194100
$ws.Range("H66").Value = 339279800
$ws.Range("I66").Value = 315493400
$ws.Range("J66").Value = 307002100

$ws.Range("D72").Value = 42586500
$ws.Range("E72").Value = 42278600
$ws.Range("F72").Value = 40447500
$ws.Range("G72").Value = 39709300
$ws.Range("H72").Value = 39609400
$ws.Range("I72").Value = 42510000
$ws.Range("J72").Value = 41952800

$ws.Range("D76").Value = 47937600
$ws.Range("E76").Value = 47629700
$ws.Range("F76").Value = 45798600
$ws.Range("G76").Value = 45060400
$ws.Range("H76").Value = 44960500
$ws.Range("I76").Value = 47861000
$ws.Range("J76").Value = 46147600

$ws.Range("D81").Value = 2437300
$ws.Range("E81").Value = 2969800
$ws.Range("F81").Value = 1816800
$ws.Range("G81").Value = 1295900
$ws.Range("H81").Value = 1542700
$ws.Range("I81").Value = 476400
$ws.Range("J81").Value = 4673300

$ws.Range("D89").Value = 6555600
$ws.Range("E89").Value = 8622300
$ws.Range("F89").Value = 21206600
$ws.Range("G89").Value = 15179100
$ws.Range("H89").Value = -1956000
$ws.Range("I89").Value = -6494300
$ws.Range("J89").Value = 3702500

$ws.Range("D94").Value = -5976300
$ws.Range("E94").Value = 3386400
$ws.Range("F94").Value = -29161700
$ws.Range("G94").Value = -16877300
$ws.Range("H94").Value = 1338400
$ws.Range("I94").Value = 7794600
$ws.Range("J94").Value = 2155800

$ws.Range("D96").Value = -1264200
$ws.Range("E96").Value = -292700
$ws.Range("F96").Value = -770500
$ws.Range("G96").Value = -493000
$ws.Range("H96").Value = -1326000
$ws.Range("I96").Value = -469100
$ws.Range("J96").Value = -1236300

$ws.Range("D100").Value = 810600
$ws.Range("E100").Value = -11886300
$ws.Range("F100").Value = 5119600
$ws.Range("G100").Value = 913800
$ws.Range("H100").Value = 16800
$ws.Range("I100").Value = -1228500
$ws.Range("J100").Value = -1652800

$ws.Range("D101").Value = 3300
$ws.Range("E101").Value = -10800
$ws.Range("F101").Value = -15700
$ws.Range("G101").Value = 34300
$ws.Range("H101").Value = 8000
$ws.Range("I101").Value = 36100

$ws.Range("D102").Value = 1393100
$ws.Range("E102").Value = 111600
$ws.Range("F102").Value = -2851300
$ws.Range("G102").Value = -750200
$ws.Range("H102").Value = -592900
$ws.Range("I102").Value = 107900
$ws.Range("J102").Value = 4205200
